$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "R02243: [] | ACONITATE-DECARBOXYLASE-RXN_reverse: [B3IUN8, P54987, Q0C8L3] | "
$ws.Range("G2").Value = "RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | "
$ws.Range("H2").Value = "RXN-16543_reverse: [] | "

$ws.Range("F3").Value = "RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | "
$ws.Range("G3").Value = "RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | "
$ws.Range("H3").Value = "RXN-16543_reverse: [] | "
